$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "26.752.22"
$ws.Range("E2").Value = "  +0.31%  "

# Row 3
$ws.Range("D3").Value = "1.603.98"
$ws.Range("E3").Value = "  +0.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.19%  "

# Row 6
$ws.Range("E6").Value = "  +0.14%  "

# Row 7
$ws.Range("E7").Value = "  +0.27%  "

# Row 8
$ws.Range("E8").Value = "  +0.00%  "

# Row 9
$ws.Range("E9").Value = "  +0.01%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "

# Row 11
$ws.Range("E11").Value = "  +0.71%  "

# Row 12
$ws.Range("D12").Value = "1.829.06"
$ws.Range("E12").Value = "  +0.45%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.07"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.99%  "

# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.583.25"
$ws.Range("E14").Value = "  -0.52%  "

# Row 15
$ws.Range("E15").Value = "  +0.46%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.02%  "

# Row 17
$ws.Range("E17").Value = "  -0.82%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "209.34"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.45%  "

# Row 19
$ws.Range("E19").Value = "  +0.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.66%  "

# Row 21
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("E22").Value = "  -4.77%  "

# Row 23
$ws.Range("E23").Value = "  +0.61%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.41%  "

# Row 25
$ws.Range("E25").Value = "  +0.14%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.10"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.24%  "

# Row 27
$ws.Range("E27").Value = "  -0.13%  "

# Row 28
$ws.Range("E28").Value = "  +0.22%  "

# Row 29
$ws.Range("E29").Value = "  -1.59%  "

# Row 30
$ws.Range("E30").Value = "  +0.17%  "

# Row 31
$ws.Range("E31").Value = "  +1.10%  "

# Row 32
$ws.Range("E32").Value = "  +0.47%  "

# Row 33
$ws.Range("D33").Value = "1.289.26"
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("E34").Value = "  +1.38%  "

# Row 35
$ws.Range("E35").Value = "  +16.27%  "

# Row 36
$ws.Range("E36").Value = "  +0.38%  "

# Row 37
$ws.Range("E37").Value = "  -5.23%  "

# Row 38
$ws.Range("E38").Value = "  -0.65%  "

# Row 39
$ws.Range("E39").Value = "  -0.36%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.44"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.04%  "

# Row 41
$ws.Range("E41").Value = "  +0.09%  "

# Row 42
$ws.Range("E42").Value = "  -0.49%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "62.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.23%  "

# Row 44
$ws.Range("D44").Value = "1.740.85"
$ws.Range("E44").Value = "  +0.52%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "90.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.83%  "

# Row 46
$ws.Range("E46").Value = "  +0.24%  "

# Row 47
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₆0103"
$ws.Range("E47").Value = "  -3.30%  "

# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.102"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.60%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.53"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.69%  "

# Row 51
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.23%  "
